# chore: update Sheets via scheduled runner
# Refreshes cached market-price / profit figures (columns H-N) on a handful
# of leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, matching a
# fresh data pull. Some rows gain or lose a cell (M/N) entirely when a value
# becomes zero/defined, so those are cleared rather than zeroed.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1118.8966
$ws.Range("I33").Value = 1031.0869
$ws.Range("K33").Value = 1031.0869
$ws.Range("M33").Value = -802.0869

$ws.Range("H64").Value = 4933.3335
$ws.Range("J64").Value = 5000
$ws.Range("L64").Value = 5000
$ws.Range("N64").Value = -5496

$ws.Range("H67").Value = 4933.3335
$ws.Range("J67").Value = 5000
$ws.Range("L67").Value = 5000
$ws.Range("N67").Value = -6716

$ws.Range("H106").Value = 4999.25
$ws.Range("I106").Value = 2000
$ws.Range("K106").Value = 2000
$ws.Range("M106").Value = -1369

$ws.Range("H111").Value = 3236.1428
$ws.Range("I111").Value = 2204.6
$ws.Range("J111").Value = 5815
$ws.Range("K111").Value = 6613.799999999999
$ws.Range("L111").Value = 17445
$ws.Range("M111").Value = -3546.799999999999
$ws.Range("N111").Value = -23579

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 2483.3333
$ws.Range("J12").Value = 5450
$ws.Range("L12").Value = 5450
$ws.Range("N12").Value = -5796

$ws.Range("H19").Value = 1000
$ws.Range("J19").Value = 1000
$ws.Range("L19").Value = 1000
$ws.Range("N19").Value = -1458

$ws.Range("H21").Value = 5000
$ws.Range("I21").Value = 5000
$ws.Range("K21").Value = 5000
$ws.Range("M21").Value = -4626

$ws.Range("H22").Value = 3000
$ws.Range("I22").Value = 1000
$ws.Range("K22").Value = 1000
$ws.Range("M22").Value = -701

$ws.Range("H25").Value = 1783
$ws.Range("I25").Value = 2666
$ws.Range("K25").Value = 2666
$ws.Range("M25").Value = -2264

$ws.Range("H29").Value = 15994
$ws.Range("I29").Value = 15994
$ws.Range("K29").Value = 15994
$ws.Range("M29").Value = -15686

$ws.Range("H30").Value = 2278
$ws.Range("I30").Value = 1130
$ws.Range("J30").Value = 4000
$ws.Range("K30").Value = 1130
$ws.Range("L30").Value = 4000
$ws.Range("M30").Value = -980
$ws.Range("N30").Value = -4300

$ws.Range("H32").Value = 12506525
$ws.Range("I32").Value = 12827199
$ws.Range("K32").Value = 12827199
$ws.Range("M32").Value = -12826912

$ws.Range("H122").Value = 2021
$ws.Range("I122").Value = 1600
$ws.Range("K122").Value = 4800
$ws.Range("M122").Value = -2350

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 1001.5
$ws.Range("J12").Value = 1953
$ws.Range("L12").Value = 1953
$ws.Range("N12").Value = -2289

$ws.Range("H18").Value = 6500
$ws.Range("J18").Value = 6500
$ws.Range("L18").Value = 6500
$ws.Range("N18").Value = -7558

$ws.Range("H24").Value = 8000
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 8000
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 8000
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -8470

$ws.Range("H99").Value = 3126.6667
$ws.Range("I99").Value = 2240
$ws.Range("K99").Value = 2240
$ws.Range("M99").Value = -742

$ws.Range("H122").Value = 39990
$ws.Range("J122").Value = 39990
$ws.Range("L122").Value = 39990
$ws.Range("N122").Value = -49790

$ws.Range("H125").Value = 53832.668
$ws.Range("J125").Value = 53832.668
$ws.Range("L125").Value = 53832.668
$ws.Range("N125").Value = -63672.668

$ws.Range("H126").Value = 38000
$ws.Range("J126").Value = 38000
$ws.Range("L126").Value = 38000
$ws.Range("N126").Value = -47880

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 870709.75
$ws.Range("I31").Value = 11301.733
$ws.Range("K31").Value = 11301.733
$ws.Range("M31").Value = -11006.733

$ws.Range("H34").Value = 870709.75
$ws.Range("I34").Value = 11301.733
$ws.Range("K34").Value = 11301.733
$ws.Range("M34").Value = -11099.733

$ws.Range("H58").Value = 1860.5
$ws.Range("I58").Value = 1499.3334
$ws.Range("K58").Value = 1499.3334
$ws.Range("M58").Value = -1296.3334

$ws.Range("H62").Value = 3199.75
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 3266.3333
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 3266.3333
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -4514.3333

$ws.Range("H65").Value = 3199.75
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 3266.3333
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 16331.6665
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -22571.6665

$ws.Range("H132").Value = 2267.2942
$ws.Range("I132").Value = 2221.5
$ws.Range("K132").Value = 6664.5
$ws.Range("M132").Value = -4134.5

$ws.Range("H134").Value = 559364.3
$ws.Range("I134").Value = 834796.5
$ws.Range("J134").Value = 8500
$ws.Range("K134").Value = 2504389.5
$ws.Range("L134").Value = 25500
$ws.Range("M134").Value = -2501854.5
$ws.Range("N134").Value = -30570

$ws.Range("H136").Value = 1860.5
$ws.Range("I136").Value = 1499.3334
$ws.Range("K136").Value = 4498.0002
$ws.Range("M136").Value = -1948.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 900
$ws.Range("I23").Value = 850
$ws.Range("K23").Value = 2550
$ws.Range("M23").Value = -2315

$ws.Range("H104").Value = 4499.6665
$ws.Range("I104").Value = 4499.6665
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 13498.9995
$ws.Range("L104").Value = 0
$ws.Range("M104").Value = -10877.9995
$ws.Range("N104").ClearContents()

$ws.Range("H105").Value = 10000
$ws.Range("J105").Value = 10000
$ws.Range("L105").Value = 30000
$ws.Range("N105").Value = -35242

$ws.Range("H140").Value = 2133.6
$ws.Range("I140").Value = 2133.6
$ws.Range("K140").Value = 6400.799999999999
$ws.Range("M140").Value = -1220.799999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 30000.5
$ws.Range("J5").Value = 40000
$ws.Range("L5").Value = 40000
$ws.Range("N5").Value = -40224

$ws.Range("H43").Value = 34000
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

$ws.Range("H80").Value = 18792.5
$ws.Range("I80").Value = 12366.111
$ws.Range("J80").Value = 30360
$ws.Range("K80").Value = 12366.111
$ws.Range("L80").Value = 30360
$ws.Range("M80").Value = -11368.111
$ws.Range("N80").Value = -32356

$ws.Range("H83").Value = 18792.5
$ws.Range("I83").Value = 12366.111
$ws.Range("J83").Value = 30360
$ws.Range("K83").Value = 61830.55500000001
$ws.Range("L83").Value = 151800
$ws.Range("M83").Value = -56838.55500000001
$ws.Range("N83").Value = -161784

$ws.Range("H113").Value = 4500
$ws.Range("I113").Value = 4500
$ws.Range("K113").Value = 4500
$ws.Range("M113").Value = -2330

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()

$ws.Range("H43").Value = 2895711.5
$ws.Range("J43").Value = 35980
$ws.Range("L43").Value = 35980
$ws.Range("N43").Value = -36366

$ws.Range("H46").Value = 2202.121
$ws.Range("I46").Value = 1585.3846
$ws.Range("K46").Value = 1585.3846
$ws.Range("M46").Value = -1397.3846

$ws.Range("H68").Value = 3474.1765
$ws.Range("I68").Value = 3048.5833
$ws.Range("J68").Value = 4495.6
$ws.Range("K68").Value = 3048.5833
$ws.Range("L68").Value = 4495.6
$ws.Range("M68").Value = -2299.5833
$ws.Range("N68").Value = -5993.6

$ws.Range("H71").Value = 3474.1765
$ws.Range("I71").Value = 3048.5833
$ws.Range("J71").Value = 4495.6
$ws.Range("K71").Value = 15242.9165
$ws.Range("L71").Value = 22478
$ws.Range("M71").Value = -11498.9165
$ws.Range("N71").Value = -29966

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 8534
$ws.Range("J132").Value = 26725.25
$ws.Range("L132").Value = 80175.75
$ws.Range("N132").Value = -85235.75
